# Updated for latest meeting: add a new "15-mars" attendance column (E)
# to the Tableau1 table and fill in attendance values for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one column (B3:D14 -> B3:E14). This updates the table
# definition (tableColumns/autoFilter/filterColumn) automatically.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:E14"))

# Set the new column header text (also renames the table column) and give
# it the same date number format as the existing "01-mars" column (D).
$ws.Range("E3").Value = "15-mars"
$ws.Range("E3").NumberFormat = $ws.Range("D3").NumberFormat

# Fill in the attendance values for the new "15-mars" date.
$ws.Range("E4").Value = 0      # Alexandre
$ws.Range("E5").Value = 1      # Robert
$ws.Range("E6").Value = 0.5    # Olivia
$ws.Range("E7").Value = 1      # Beenita
$ws.Range("E8").Value = 1      # Anushan
$ws.Range("E9").Value = 1      # Vytas
$ws.Range("E10").Value = 1     # Juliano
$ws.Range("E11").Value = 1     # Jakub
$ws.Range("E12").Value = 1     # Mitchell
$ws.Range("E14").Value = 1     # TA

# Extend the colour-scale conditional formatting over the new column.
$fc = $ws.Range("C4:D14").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("C4:E14"))

# Match the selection/active cell recorded after the edit.
[void]$ws.Range("E6").Select()
